# Update column G (K) values on Sheet1 for rows 2-77.
# New values taken from the target OOXML diff (old Strike# values replaced
# with the regenerated K values).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newK = @{
    2  = 0
    3  = 1
    4  = 2
    5  = 0
    6  = 0
    7  = 0
    8  = 1
    9  = 1
    10 = 0
    11 = 0
    12 = 2
    13 = 3
    14 = 0
    15 = 0
    16 = 2
    17 = 1
    18 = 3
    19 = 2
    20 = 1
    21 = 1
    22 = 0
    23 = 1
    24 = 1
    25 = 0
    26 = 3
    27 = 1
    28 = 0
    29 = 2
    30 = 2
    31 = 2
    32 = 1
    33 = 3
    34 = 0
    35 = 2
    36 = 1
    37 = 2
    38 = 0
    39 = 1
    40 = 1
    41 = 1
    42 = 1
    43 = 1
    44 = 0
    45 = 1
    46 = 1
    47 = 1
    48 = 3
    49 = 1
    50 = 3
    51 = 4
    52 = 1
    53 = 2
    54 = 1
    55 = 4
    56 = 4
    57 = 4
    58 = 3
    59 = 2
    60 = 0
    61 = 4
    62 = 2
    63 = 1
    64 = 1
    65 = 1
    66 = 1
    67 = 0
    68 = 2
    69 = 1
    70 = 2
    71 = 2
    72 = 1
    73 = 2
    74 = 2
    75 = 1
    76 = 0
    77 = 0
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
